# Scheduled market-price refresh for Bahamut_Profits sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) with the latest
# Universalis price snapshot for the affected leve rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 560.125
$ws.Range("I18").Value = 576.2
$ws.Range("J18").Value = 533.3333
$ws.Range("K18").Value = 576.2
$ws.Range("L18").Value = 533.3333
$ws.Range("M18").Value = -292.2
$ws.Range("N18").Value = -1101.3333
$ws.Range("H69").Value = 4559.8
$ws.Range("I69").Value = 3999.5
$ws.Range("J69").Value = 4933.3335
$ws.Range("K69").Value = 11998.5
$ws.Range("L69").Value = 14800.0005
$ws.Range("M69").Value = -11124.5
$ws.Range("N69").Value = -16548.0005
$ws.Range("H72").Value = 4559.8
$ws.Range("I72").Value = 3999.5
$ws.Range("J72").Value = 4933.3335
$ws.Range("K72").Value = 35995.5
$ws.Range("L72").Value = 44400.0015
$ws.Range("M72").Value = -31627.5
$ws.Range("N72").Value = -53136.0015
$ws.Range("H74").Value = 3464521.5
$ws.Range("I74").Value = 3996952.5
$ws.Range("J74").Value = 3719.5
$ws.Range("K74").Value = 3996952.5
$ws.Range("L74").Value = 3719.5
$ws.Range("M74").Value = -3996016.5
$ws.Range("N74").Value = -5591.5
$ws.Range("H77").Value = 3464521.5
$ws.Range("I77").Value = 3996952.5
$ws.Range("J77").Value = 3719.5
$ws.Range("K77").Value = 19984762.5
$ws.Range("L77").Value = 18597.5
$ws.Range("M77").Value = -19980082.5
$ws.Range("N77").Value = -27957.5
$ws.Range("H116").Value = 4032.0454
$ws.Range("I116").Value = 3953.2666
$ws.Range("K116").Value = 3953.2666
$ws.Range("M116").Value = -511.2665999999999
$ws.Range("H137").Value = 848.26666
$ws.Range("I137").Value = 823.1429000000001
$ws.Range("K137").Value = 2469.4287
$ws.Range("M137").Value = 80.57129999999961

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 19419.6
$ws.Range("J44").Value = 19419.6
$ws.Range("L44").Value = 19419.6
$ws.Range("N44").Value = -20395.6
$ws.Range("H61").Value = 1338.4783
$ws.Range("I61").Value = 939.4706
$ws.Range("J61").Value = 2469
$ws.Range("K61").Value = 939.4706
$ws.Range("L61").Value = 2469
$ws.Range("M61").Value = -727.4706
$ws.Range("N61").Value = -2893
$ws.Range("H74").Value = 1308.5172
$ws.Range("I74").Value = 1334.125
$ws.Range("J74").Value = 1185.6
$ws.Range("K74").Value = 1334.125
$ws.Range("L74").Value = 1185.6
$ws.Range("M74").Value = -460.125
$ws.Range("N74").Value = -2933.6
$ws.Range("H77").Value = 1308.5172
$ws.Range("I77").Value = 1334.125
$ws.Range("J77").Value = 1185.6
$ws.Range("K77").Value = 6670.625
$ws.Range("L77").Value = 5928
$ws.Range("M77").Value = -2302.625
$ws.Range("N77").Value = -14664
$ws.Range("H97").Value = 1063
$ws.Range("I97").Value = 1146
$ws.Range("J97").Value = 980
$ws.Range("K97").Value = 1146
$ws.Range("L97").Value = 980
$ws.Range("M97").Value = -650
$ws.Range("N97").Value = -1972
$ws.Range("H102").Value = 1980.7693
$ws.Range("I102").Value = 1644.7894
$ws.Range("J102").Value = 2892.7144
$ws.Range("K102").Value = 1644.7894
$ws.Range("L102").Value = 2892.7144
$ws.Range("M102").Value = -22.78939999999989
$ws.Range("N102").Value = -6136.7144
$ws.Range("H122").Value = 821.8570999999999
$ws.Range("I122").Value = 847.8182
$ws.Range("J122").Value = 726.6667
$ws.Range("K122").Value = 2543.4546
$ws.Range("L122").Value = 2180.0001
$ws.Range("M122").Value = -93.45460000000003
$ws.Range("N122").Value = -7080.0001
$ws.Range("H136").Value = 1338.4783
$ws.Range("I136").Value = 939.4706
$ws.Range("J136").Value = 2469
$ws.Range("K136").Value = 2818.4118
$ws.Range("L136").Value = 7407
$ws.Range("M136").Value = -268.4117999999999
$ws.Range("N136").Value = -12507

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22434.715
$ws.Range("I134").Value = 1686.8948
$ws.Range("J134").Value = 94109
$ws.Range("K134").Value = 5060.6844
$ws.Range("L134").Value = 282327
$ws.Range("M134").Value = -2525.6844
$ws.Range("N134").Value = -287397

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2679.5483
$ws.Range("I31").Value = 2794.2144
$ws.Range("K31").Value = 2794.2144
$ws.Range("M31").Value = -2499.2144
$ws.Range("H34").Value = 2679.5483
$ws.Range("I34").Value = 2794.2144
$ws.Range("K34").Value = 2794.2144
$ws.Range("M34").Value = -2592.2144
$ws.Range("H35").Value = 2043.75
$ws.Range("I35").Value = 425
$ws.Range("J35").Value = 6900
$ws.Range("K35").Value = 425
$ws.Range("L35").Value = 6900
$ws.Range("M35").Value = -131
$ws.Range("N35").Value = -7488
$ws.Range("H50").Value = 9343
$ws.Range("J50").Value = 9343
$ws.Range("L50").Value = 9343
$ws.Range("N50").Value = -10593
$ws.Range("H51").Value = 8557.200000000001
$ws.Range("I51").Value = 7295
$ws.Range("J51").Value = 9398.666999999999
$ws.Range("K51").Value = 7295
$ws.Range("L51").Value = 9398.666999999999
$ws.Range("M51").Value = -6559
$ws.Range("N51").Value = -10870.667
$ws.Range("H60").Value = 5984.1665
$ws.Range("I60").Value = 3533.3333
$ws.Range("J60").Value = 8435
$ws.Range("K60").Value = 3533.3333
$ws.Range("L60").Value = 8435
$ws.Range("M60").Value = -3022.3333
$ws.Range("N60").Value = -9457
$ws.Range("H61").Value = 8557.200000000001
$ws.Range("I61").Value = 7295
$ws.Range("J61").Value = 9398.666999999999
$ws.Range("K61").Value = 7295
$ws.Range("L61").Value = 9398.666999999999
$ws.Range("M61").Value = -6947
$ws.Range("N61").Value = -10094.667
$ws.Range("H134").Value = 967.53625
$ws.Range("I134").Value = 827.56604
$ws.Range("J134").Value = 1431.1875
$ws.Range("K134").Value = 2482.69812
$ws.Range("L134").Value = 4293.5625
$ws.Range("M134").Value = 52.30187999999998
$ws.Range("N134").Value = -9363.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 792.62
$ws.Range("I131").Value = 352.2
$ws.Range("J131").Value = 870.3412
$ws.Range("K131").Value = 1056.6
$ws.Range("L131").Value = 2611.0236
$ws.Range("M131").Value = 3983.4
$ws.Range("N131").Value = -12691.0236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3500.1875
$ws.Range("I80").Value = 3500.2144
$ws.Range("K80").Value = 3500.2144
$ws.Range("M80").Value = -2502.2144
$ws.Range("H83").Value = 3500.1875
$ws.Range("I83").Value = 3500.2144
$ws.Range("K83").Value = 17501.072
$ws.Range("M83").Value = -12509.072
$ws.Range("H97").Value = 1676.1765
$ws.Range("I97").Value = 1705.5
$ws.Range("J97").Value = 1634.2858
$ws.Range("K97").Value = 1705.5
$ws.Range("L97").Value = 1634.2858
$ws.Range("M97").Value = -1209.5
$ws.Range("N97").Value = -2626.2858
$ws.Range("H122").Value = 8115116.5
$ws.Range("I122").Value = 9979045
$ws.Range("J122").Value = 6251187.5
$ws.Range("K122").Value = 29937135
$ws.Range("L122").Value = 18753562.5
$ws.Range("M122").Value = -29934685
$ws.Range("N122").Value = -18758462.5
$ws.Range("H123").Value = 11269.833
$ws.Range("J123").Value = 11269.833
$ws.Range("L123").Value = 11269.833
$ws.Range("N123").Value = -16169.833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7368.9473
$ws.Range("I122").Value = 9200.714
$ws.Range("J122").Value = 2240
$ws.Range("K122").Value = 27602.142
$ws.Range("L122").Value = 6720
$ws.Range("M122").Value = -25152.142
$ws.Range("N122").Value = -11620
$ws.Range("H132").Value = 2062.9187
$ws.Range("I132").Value = 1717.3062
$ws.Range("J132").Value = 2520.6216
$ws.Range("K132").Value = 5151.9186
$ws.Range("L132").Value = 7561.864799999999
$ws.Range("M132").Value = -2621.9186
$ws.Range("N132").Value = -12621.8648

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents() | Out-Null
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents() | Out-Null
$ws.Range("H132").Value = 1558.75
$ws.Range("I132").Value = 1049.2285
$ws.Range("J132").Value = 2272.08
$ws.Range("K132").Value = 3147.6855
$ws.Range("L132").Value = 6816.24
$ws.Range("M132").Value = -617.6854999999996
$ws.Range("N132").Value = -11876.24
$ws.Range("H136").Value = 2682.45
$ws.Range("I136").Value = 2707.9333
$ws.Range("J136").Value = 2606
$ws.Range("K136").Value = 8123.7999
$ws.Range("L136").Value = 8910
$ws.Range("M136").Value = -5573.7999
$ws.Range("N136").Value = -12918

Write-Host "Bahamut_Profits price refresh applied."
